# Updates the cryptocurrency price/volume table on Sheet1 with the latest
# scraped values from coinranking.com (GitHub Actions scheduled refresh).
#
# Columns: B=Coin name, C=Link, D=Price, E=Volume(1h) change.
# Most rows only refresh Price (D) and Volume (E); a few rows had their
# coin re-ranked and so also get new Coin/Link text (B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $value into column $colLetter of row $rowNum as literal text.
# Price values such as "10.50" or "1.00" would otherwise be auto-parsed by
# Excel as numbers (dropping the significant trailing zero), so the cell is
# explicitly formatted as Text before the assignment, then restored to the
# default "Normal" style so no stray formatting is left behind.
function Set-TextCell($rowNum, $colLetter, $value) {
    $cell = $ws.Range("$colLetter$rowNum")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$rowUpdates = @(
    @{ Row=2; D='66.858.59'; E='  +1.46%  ' }
    @{ Row=3; D='3.881.06'; E='  +2.00%  ' }
    @{ Row=4; D='0.998'; E='  -0.42%  ' }
    @{ Row=5; D='471.32'; E='  +10.01%  ' }
    @{ Row=6; D='145.46'; E='  +10.91%  ' }
    @{ Row=7; D='0.629'; E='  +3.10%  ' }
    @{ Row=8; D='0.998'; E='  -0.22%  ' }
    @{ Row=9; D='0.746'; E='  +1.23%  ' }
    @{ Row=10; D='0.157'; E='  +0.93%  ' }
    @{ Row=11; D='0.0000314'; E='  -4.31%  ' }
    @{ Row=12; D='43.52'; E='  +3.27%  ' }
    @{ Row=13; D='10.50'; E='  -2.23%  ' }
    @{ Row=14; D='4.481.73'; E='  +1.16%  ' }
    @{ Row=15; D='15.03'; E='  -1.91%  ' }
    @{ Row=16; D='3.893.96'; E='  +2.23%  ' }
    @{ Row=17; E='  -0.23%  ' }
    @{ Row=18; D='20.18'; E='  +0.52%  ' }
    @{ Row=19; E='  +3.50%  ' }
    @{ Row=20; D='67.168.87'; E='  +1.48%  ' }
    @{ Row=21; D='435.26'; E='  +5.46%  ' }
    @{ Row=22; D='15.10'; E='  -3.36%  ' }
    @{ Row=23; D='3.36'; E='  +5.31%  ' }
    @{ Row=24; D='88.22'; E='  +2.83%  ' }
    @{ Row=25; D='3.58'; E='  +8.10%  ' }
    @{ Row=26; D='38.62'; E='  +4.34%  ' }
    @{ Row=27; D='10.12'; E='  +7.13%  ' }
    @{ Row=28; D='10.12'; E='  +1.51%  ' }
    @{ Row=29; D='5.55'; E='  +2.70%  ' }
    @{ Row=30; D='733.33'; E='  +2.61%  ' }
    @{ Row=31; D='13.97'; E='  -2.82%  ' }
    @{ Row=32; E='  +6.77%  ' }
    @{ Row=33; E='  +0.39%  ' }
    @{ Row=34; D='43.53'; E='  +11.40%  ' }
    @{ Row=35; D='0.162'; E='  +5.63%  ' }
    @{ Row=36; D='59.20'; E='  +5.33%  ' }
    @{ Row=37; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='5.56'; E='  -8.09%  ' }
    @{ Row=38; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  +0.12%  ' }
    @{ Row=39; D='0.0486'; E='  +2.89%  ' }
    @{ Row=40; D='3.03'; E='  +3.58%  ' }
    @{ Row=41; D='0.348'; E='  +7.04%  ' }
    @{ Row=42; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0₃0694'; E='  -0.31%  ' }
    @{ Row=43; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='2.61'; E='  +3.92%  ' }
    @{ Row=44; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.142'; E='  +3.12%  ' }
    @{ Row=45; E='  -0.21%  ' }
    @{ Row=46; D='3.49'; E='  +2.55%  ' }
    @{ Row=47; D='2.21'; E='  +7.41%  ' }
    @{ Row=48; D='2.77'; E='  +5.00%  ' }
    @{ Row=49; E='  -0.82%  ' }
    @{ Row=50; D='2.93'; E='  +3.55%  ' }
    @{ Row=51; D='142.78'; E='  +0.98%  ' }
)

foreach ($update in $rowUpdates) {
    $rowNum = $update.Row

    if ($update.ContainsKey('B')) {
        $ws.Range("B$rowNum").Value = $update.B
    }
    if ($update.ContainsKey('C')) {
        $ws.Range("C$rowNum").Value = $update.C
    }
    if ($update.ContainsKey('D')) {
        Set-TextCell $rowNum "D" $update.D
    }
    if ($update.ContainsKey('E')) {
        $ws.Range("E$rowNum").Value = $update.E
    }
}
